# Regenerate save_data column G ("K") with newly calculated strikeout values (s_vals),
# replacing the previous "Strike#" derived figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..75 (row index -> value), in order.
$s_vals = @(
    1,1,0,2,1,0,3,2,2,1,
    3,1,0,3,3,0,3,1,0,2,
    0,1,1,1,1,2,1,3,1,4,
    1,0,0,1,2,1,2,1,0,3,
    1,1,1,1,1,0,2,2,1,1,
    1,1,3,3,0,0,0,2,4,0,
    1,2,0,2,0,1,1,0,1,0,
    2,1,0,1
)

$startRow = 2
for ($i = 0; $i -lt $s_vals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $s_vals[$i]
}
